# Actualización de horarios Línea 141 - 1267
# Scrape timestamp: 02:38:35 -> 03:27:48

$wb = $excel.ActiveWorkbook

$nuevaHora = "03:27:48"

# ---------------------------------------------------------------------------
# Hoja "LP1912": nuevos datos de arribos + actualización de encabezado
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $nuevaHora"
$ws1.Range("A3").Value = "Total filas: 6"

$filas = @(
    @("03:27:48", "03:48", "14_ABASTO",      21, "LP1912"),
    @("03:27:48", "04:01", "81_EL PELIGRO",  34, "LP1912"),
    @("03:27:48", "04:47", "81_EL PELIGRO",  80, "LP1912"),
    @("03:27:48", "04:53", "11_ETCHEVERRY",  86, "LP1912"),
    @("03:27:48", "05:16", "17_ROMERO",     109, "LP1912"),
    @("03:27:48", "05:22", "23_HERNANDEZ",  115, "LP1912")
)

$r = 6
foreach ($fila in $filas) {
    $ws1.Cells.Item($r, 1).Value = $fila[0]
    $ws1.Cells.Item($r, 2).Value = $fila[1]
    $ws1.Cells.Item($r, 3).Value = $fila[2]
    $ws1.Cells.Item($r, 4).Value = $fila[3]
    $ws1.Cells.Item($r, 5).Value = $fila[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Hoja "LP1912-215": solo se actualiza la hora del scrap
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $nuevaHora"

# ---------------------------------------------------------------------------
# Hoja "6203-6173": solo se actualiza la hora del scrap
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $nuevaHora"
